# Insert a new row at position 107, shifting existing rows 107..165 down to 108..166.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(107).Insert()

# Populate the newly inserted row 107 with its data.
$ws.Cells.Item(107, 1).Value = 4
$ws.Cells.Item(107, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(107, 3).Value = "Los Lagos"
$ws.Cells.Item(107, 4).Value = 44523
$ws.Cells.Item(107, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(107, 5).Value = 10
$ws.Cells.Item(107, 6).Value = "Fruta"
$ws.Cells.Item(107, 7).Value = 100108
$ws.Cells.Item(107, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(107, 9).Value = 100108005
$ws.Cells.Item(107, 10).Value = "Piña"
$ws.Cells.Item(107, 11).Value = "Caramelo"
$ws.Cells.Item(107, 12).Value = "Tercera"
$ws.Cells.Item(107, 13).Value = 300
$ws.Cells.Item(107, 14).Value = 20000
$ws.Cells.Item(107, 15).Value = 21000
$ws.Cells.Item(107, 16).Value = 20500
$ws.Cells.Item(107, 17).Value = "`$/caja 16 unidades"
$ws.Cells.Item(107, 18).Value = "Ecuador"
$ws.Cells.Item(107, 19).Value = 1281
$ws.Cells.Item(107, 20).Value = 16
